# Reorder the MSME indicator rows on the Summary sheet so that
# "Enterprises density (per 1000 people)" / "2.4" moves up to
# immediately follow the Micro/SMEs/MSMEs header row, pushing
# "Employment (% of total)" and "Enterprises (absolute #)" down
# by one row. "Enterprises (% of total)" / "99.4" stays in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values of the affected rows (A and D columns) before
# overwriting anything. Use Value2 since plain Value does not reliably
# resolve to a scalar in this COM runtime.
$density_label = $ws.Range("A14").Value2   # "Enterprises density (per 1000 people)"
$density_value = $ws.Range("D14").Value2   # "2.4"

$employment_label = $ws.Range("A12").Value2   # "Employment (% of total)"
$employment_value = $ws.Range("D12").Value2   # "80"

$enterprises_label = $ws.Range("A13").Value2  # "Enterprises (absolute #)"
$enterprises_value = $ws.Range("D13").Value2  # "126237"

# Write the new order: density moves to row 12, employment moves to
# row 13, enterprises (absolute #) moves to row 14. Row 15 is untouched.
# The D-column values are stored as text (not numbers) in the original
# file, so prefix them with an apostrophe to force a text entry and
# keep the same cell type.
$ws.Range("A12").Value2 = $density_label
$ws.Range("D12").Value2 = "'" + $density_value

$ws.Range("A13").Value2 = $employment_label
$ws.Range("D13").Value2 = "'" + $employment_value

$ws.Range("A14").Value2 = $enterprises_label
$ws.Range("D14").Value2 = "'" + $enterprises_value
